$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $r = $ws.Range($cell)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.ClearFormats()
}

Set-TextValue 'D2' '30.406.11'
Set-TextValue 'E2' '  +0.21%  '
Set-TextValue 'D3' '1.937.23'
Set-TextValue 'E3' '  +0.18%  '
Set-TextValue 'E4' '  +0.33%  '
Set-TextValue 'D5' '0.7455'
Set-TextValue 'E5' '  +3.65%  '
Set-TextValue 'D6' '245.26'
Set-TextValue 'E6' '  -2.65%  '
Set-TextValue 'D7' '1.004'
Set-TextValue 'E7' '  +0.26%  '
Set-TextValue 'B8' 'Cardano'
Set-TextValue 'C8' 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
Set-TextValue 'D8' '0.3170'
Set-TextValue 'E8' '  -3.35%  '
Set-TextValue 'B9' 'Solana'
Set-TextValue 'C9' 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-TextValue 'D9' '27.53'
Set-TextValue 'E9' '  -0.85%  '
Set-TextValue 'D10' '0.06986'
Set-TextValue 'E10' '  -3.91%  '
Set-TextValue 'E11' '  -3.15%  '
Set-TextValue 'D12' '0.07998'
Set-TextValue 'E12' '  -1.14%  '
Set-TextValue 'D13' '1.937.80'
Set-TextValue 'E13' '  +0.17%  '
Set-TextValue 'D14' '5.360'
Set-TextValue 'E14' '  -1.30%  '
Set-TextValue 'D15' '94.47'
Set-TextValue 'E15' '  -0.15%  '
Set-TextValue 'D16' '14.42'
Set-TextValue 'E16' '  -4.19%  '
Set-TextValue 'D17' '30.391.75'
Set-TextValue 'E17' '  +0.18%  '
Set-TextValue 'D18' '252.82'
Set-TextValue 'E18' '  -0.19%  '
Set-TextValue 'D19' '0.000007937'
Set-TextValue 'E19' '  -3.54%  '
Set-TextValue 'D20' '5.733'
Set-TextValue 'E20' '  -1.30%  '
Set-TextValue 'D21' '2.190.72'
Set-TextValue 'E21' '  +0.24%  '
Set-TextValue 'E22' '  +0.19%  '
Set-TextValue 'D23' '1.003'
Set-TextValue 'E23' '  +0.29%  '
Set-TextValue 'D24' '6.674'
Set-TextValue 'E24' '  -3.94%  '
Set-TextValue 'D25' '9.498'
Set-TextValue 'E25' '  -2.32%  '
Set-TextValue 'D26' '166.11'
Set-TextValue 'E26' '  -0.18%  '
Set-TextValue 'D27' '18.98'
Set-TextValue 'E27' '  -1.53%  '
Set-TextValue 'D28' '0.1324'
Set-TextValue 'E28' '  +2.42%  '
Set-TextValue 'D29' '2.239'
Set-TextValue 'E29' '  -4.48%  '
Set-TextValue 'D30' '1.364'
Set-TextValue 'E30' '  +0.28%  '
Set-TextValue 'D31' '1.510'
Set-TextValue 'E31' '  -2.29%  '
Set-TextValue 'D32' '4.359'
Set-TextValue 'E32' '  -1.87%  '
Set-TextValue 'D33' '4.109'
Set-TextValue 'E33' '  -2.25%  '
Set-TextValue 'D34' '0.05154'
Set-TextValue 'E34' '  -1.47%  '
Set-TextValue 'D35' '1.271'
Set-TextValue 'E35' '  +0.09%  '
Set-TextValue 'D36' '0.7463'
Set-TextValue 'E36' '  -0.39%  '
Set-TextValue 'E37' '  +0.93%  '
Set-TextValue 'D38' '0.01948'
Set-TextValue 'E38' '  -0.94%  '
Set-TextValue 'D39' '2.805'
Set-TextValue 'E39' '  +0.13%  '
Set-TextValue 'D40' '77.80'
Set-TextValue 'E40' '  -1.76%  '
Set-TextValue 'E41' '  -0.58%  '
Set-TextValue 'D42' '0.4458'
Set-TextValue 'E42' '  -2.06%  '
Set-TextValue 'E43' '  -3.55%  '
Set-TextValue 'E44' '  +0.26%  '
Set-TextValue 'D45' '0.8318'
Set-TextValue 'E45' '  -1.35%  '
Set-TextValue 'E46' '  -0.61%  '
Set-TextValue 'D47' '9.735'
Set-TextValue 'E47' '  -0.78%  '
Set-TextValue 'D48' '7.453'
Set-TextValue 'E48' '  -0.03%  '
Set-TextValue 'D49' '988.54'
Set-TextValue 'E49' '  +11.51%  '
Set-TextValue 'D50' '37.28'
Set-TextValue 'E50' '  +1.30%  '
Set-TextValue 'D51' '0.06025'
Set-TextValue 'E51' '  -0.50%  '
